$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "U2FsdGVkX1/uaCilOc9bb4XC5m5Y6DQN8bfZxHVA8sJ10SIoCpbpYjkD/3CmdJ46xhtaSBZv7CuQwilp2QwRMP7CyrHzY1U616pbYkm82+aoAa6sdI69qDyrtd7zCl7gZIt+NAniqZVzpa31YfpTZ7t99v+ziPwzc/Z00is155mUQJId4enTVgK76dck9LmoyK+bWxLYLpSovqZgZ4f7idlk7UkJBXKaVlR8wAdtOdM7YToEpstJc5mOyXF87tQyq7XIcVLaNT9LlDJoTgJL34lD8bnqdzQMZ/8xS4wExQK2PJ1Cw2WTBfJXKSKio2CFOfGM7RsWt96LA+XGRVZfE0+SgvfnrbXhCpTiWvIDe94HDOXNrU6nHfGYnOssmAExU7MJZlnumJU0oG/aiXInemQSjwtmGqAiLqkOR4FcyfctNzy6VkEBN+Fj60CLvMHg9EZlomjwupZcUrou+EjQ/Gg+fbvgLw6crqYs55hE4sftk1z4CwVp/HgFiFcnDlUgKUrhvyB6NQb6JS8u6n/ui2me0q0SIUH5Tp7NrfzQq/XFxkfOF756nBEp9lKmrBb/aD6MdSK27O+zlAbaq6ejU5Q33W5IJ/uEIQCexz/m1d5F814OhXL6osdG/BBV5HndVi5VEs3C8ON4/2ei4lq+gi7Wkqfza4Wi1kp0/YkGfGyQqkqwpvwIpo5jB5tjrIpBak5M3ONmp4zz92jR9H5heHCv28RCBa0ggy9oJgV2B0F8VtS+hvsGGNvjnY8SBOHwTjQqLZIAVs4l0SNd6ToEIVldiTYU9I015k277nQEVjgq6V+M6WlARzvNqesVRGrEl/QxChjIjaWRuJGWeQ7j89KPIso2K0ujiYwnDLof4j4VvHmnMB3gOqZbwJLhvIY6Yu7xCx9F3afBQYzXIpcvJ/M+3mX3KvVj3RXoVILgiblVKOaRjNaKveRLHFcxwwiA29AdZsngI0wmLqwb4IsDzHK0p3nyFaStFDoZTaNc2d4UtXOpTObha+WaeqQgLgogV/bQVp5UTH8hxfNoUloPjHFjkWhigGSF4gKIfiTpnSBjaOu0dGt0rZYdZvxJblmvTwi731SColrLiFCNOmE4x7gt8kqSMl5C3Ay2Hb1Y3WU6ACA05lzwX5MOqvG9kwmH2GxN516MjBjn5yc7tPWyKHg1zl5eE/fjv5WxttZthTHfvaAtpeEGQv1xMPs6DQPYI5GROMAxAfnoTnJGNgcEXFWSy7UhIdKv+Pkkgf2yvZuk3xoH4NDhk0iuyJtIdLPuYv5w/jblH0AUILavcqnrLa0e7UfxFFLjUamJg9pUbLo="

$ws.Range("B2").Value = $newValue
$ws.Range("B4").Value = $newValue
$ws.Range("B6").Value = $newValue

$ws.Range("B11").Select()
